$wb = $excel.ActiveWorkbook

function Set-TextValue($ws, $addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
}

function Set-NumValue($ws, $addr, $val) {
    $ws.Range($addr).Value = $val
}

# ------------------------------------------------------------------
# Class10 sheet
# ------------------------------------------------------------------
$ws10 = $wb.Worksheets.Item("Class10")

Set-TextValue $ws10 "C2" "11-01-1997"
Set-TextValue $ws10 "G2" "43"
Set-NumValue  $ws10 "K2" 324

Set-TextValue $ws10 "C3" "11-01-1998"
Set-NumValue  $ws10 "K3" 373

Set-TextValue $ws10 "A4" "Stephan"
Set-TextValue $ws10 "B4" "98"
Set-TextValue $ws10 "C4" "08-09-1996"
Set-TextValue $ws10 "D4" "10"
Set-TextValue $ws10 "E4" "France"
Set-TextValue $ws10 "F4" "55"
Set-TextValue $ws10 "G4" "55"
Set-TextValue $ws10 "H4" "55"
Set-TextValue $ws10 "I4" "55"
Set-TextValue $ws10 "J4" "55"
Set-NumValue  $ws10 "K4" 275
Set-TextValue $ws10 "L4" "F"

# ------------------------------------------------------------------
# Class4 sheet
# ------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("Class4")

Set-TextValue $ws4 "C2" "10-09-1991"
Set-TextValue $ws4 "F2" "100"
Set-NumValue  $ws4 "K2" 381
Set-TextValue $ws4 "L2" "C"

Set-TextValue $ws4 "C3" "10-05-1994"

Set-TextValue $ws4 "A4" "Jofin"
Set-TextValue $ws4 "B4" "15"
Set-TextValue $ws4 "C4" "10-09-2022"
Set-TextValue $ws4 "D4" "4"
Set-TextValue $ws4 "E4" "Pulickal"
Set-TextValue $ws4 "F4" "68"
Set-TextValue $ws4 "G4" "78"
Set-TextValue $ws4 "H4" "78"
Set-TextValue $ws4 "I4" "100"
Set-TextValue $ws4 "J4" "88"
Set-NumValue  $ws4 "K4" 412
Set-TextValue $ws4 "L4" "B"

# ------------------------------------------------------------------
# Class5 sheet
# ------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item("Class5")

Set-TextValue $ws5 "C2" "02-04-1990`n"
Set-NumValue  $ws5 "K2" 312

Set-TextValue $ws5 "C3" "22-08-1999`n"
Set-NumValue  $ws5 "K3" 103

Set-TextValue $ws5 "A4" "Joseph"
Set-TextValue $ws5 "B4" "56"
Set-TextValue $ws5 "C4" "12-08-1994"
Set-TextValue $ws5 "D4" "5"
Set-TextValue $ws5 "E4" "Munich"
Set-TextValue $ws5 "F4" "76"
Set-TextValue $ws5 "G4" "88"
Set-TextValue $ws5 "H4" "99"
Set-TextValue $ws5 "I4" "66"
Set-TextValue $ws5 "J4" "99"
Set-NumValue  $ws5 "K4" 428
Set-TextValue $ws5 "L4" "B"
